$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 and 20: TRON/Uniswap swap positions
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'11.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.92%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.114"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "

# Row 41 and 42: TheGraph/dogwifhat swap positions
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.94%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "'0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.51%  "

# Remaining price/volume updates
$ws.Range("D2").Value = "69.004.66"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.807.76"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'600.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'163.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").Value = "3.804.04"
$ws.Range("E7").Value = "  +1.75%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.534"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").Value = "'6.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'37.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "'0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "4.443.90"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "3.824.12"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "69.143.72"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'7.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D21").Value = "'17.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'485.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  +5.88%  "
$ws.Range("D25").Value = "'84.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -2.71%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D30").Value = "'2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").Value = "'8.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("D32").Value = "'2.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.71%  "
$ws.Range("D33").Value = "3.963.79"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").Value = "'31.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "3.752.05"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").Value = "'1.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").Value = "'0.140"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.81%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D43").Value = "'435.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").Value = "'48.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D47").Value = "'8.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "2.827.88"
$ws.Range("D49").Value = "'141.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'39.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.53%  "
$ws.Range("E51").Value = "  -0.13%  "
